$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.274.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.493.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.090.15"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.494.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.79"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.311.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.79"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.75"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.84"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.567"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.634.20"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.516.79"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.49"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.21"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.87"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.43"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0781"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.81"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.481.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.77"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.898"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.37%  "
